# Insert a new "Peer Review" section header row above the existing
# peer-review detail row (which currently carries the "Peer-review: ..."
# prefix), shifting the remaining rows down by one, and strip the
# "Peer-review: " prefix from the detail row's text since it now sits
# under its own header.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row before row 30 (the old "Peer-review: ..." row),
# pushing rows 30-37 down to 31-38.
$ws.Rows.Item(30).Insert()

# Former row 30 content (now row 31): drop the "Peer-review: " prefix.
# (Edit this cell first so the shared-string table keeps this text at its
# original slot and appends the brand-new "Peer Review" string after it.)
$ws.Range("A31").Value = "Epidemiology, European Journal of Epidemiology, Journal of Causal Inference,  BMC Medical Research Methodology, Plos One, The American Journal of Drug and Alcohol Abuse, Journal of the Intensive Care Society"

# New header row 30: "Peer Review" section heading, column D = "O"
$ws.Range("A30").Value = "Peer Review"
$ws.Range("D30").Value = "O"

# Update the selection to match the saved state (D30 is the active cell).
$ws.Range("D30").Select()
